# Daily attendance processing - reorder the "Recorded By" (column G) entries
# so that "System" / "system" sorts before the other recorder names.
#
# Two exact substitutions are applied, wherever they occur in column G:
#   "dnasr281@gmail.com, System"            -> "System, dnasr281@gmail.com"
#   "System, backup@backdoor.com, system"   -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
